$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D: Price, E: Volume) keep their text formatting
# so Excel does not auto-convert numeric-looking strings into numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

# Apply updated cell values from the refreshed crypto data feed
$ws.Range('D2').Value = '28.544.06'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '1.962.46'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').Value = '1.013'
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('D5').Value = '322.84'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '1.011'
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('D7').Value = '0.4803'
$ws.Range('E7').Value = '  -3.81%  '
$ws.Range('D8').Value = '0.4071'
$ws.Range('E8').Value = '  -3.41%  '
$ws.Range('D9').Value = '54.18'
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').Value = '0.08486'
$ws.Range('E10').Value = '  -7.91%  '
$ws.Range('E11').Value = '  -3.36%  '
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.961.15'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '7.590'
$ws.Range('E14').Value = '  -3.68%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '6.186'
$ws.Range('E15').Value = '  -3.98%  '
$ws.Range('D16').Value = '1.013'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').Value = '91.05'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('E18').Value = '  -2.65%  '
$ws.Range('D19').Value = '0.06652'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').Value = '18.54'
$ws.Range('E20').Value = '  -3.59%  '
$ws.Range('D21').Value = '1.012'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D22').Value = '5.886'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').Value = '28.596.59'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').Value = '11.49'
$ws.Range('E24').Value = '  -4.44%  '
$ws.Range('D25').Value = '2.304'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '2.284.20'
$ws.Range('E26').Value = '  +2.23%  '
$ws.Range('D27').Value = '156.36'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('D29').Value = '2.185'
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').Value = '5.883'
$ws.Range('E30').Value = '  -5.17%  '
$ws.Range('D31').Value = '124.96'
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').Value = '0.9927'
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').Value = '0.09685'
$ws.Range('E33').Value = '  -1.73%  '
$ws.Range('D34').Value = '1.461'
$ws.Range('E34').Value = '  -4.78%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.702'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '5.649'
$ws.Range('E36').Value = '  -2.22%  '
$ws.Range('D37').Value = '9.125'
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('D38').Value = '0.02337'
$ws.Range('E38').Value = '  -3.61%  '
$ws.Range('D39').Value = '0.06243'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = '1.255'
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('D41').Value = '0.6246'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('D43').Value = '1.011'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Value = '0.1920'
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('D45').Value = '1.359'
$ws.Range('E45').Value = '  +5.42%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.5968'
$ws.Range('E46').Value = '  -4.07%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '13.05'
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('D48').Value = '2.068'
$ws.Range('E48').Value = '  -5.26%  '
$ws.Range('D49').Value = '3.416'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('D50').Value = '0.06848'
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('D51').Value = '111.37'
$ws.Range('E51').Value = '  -1.07%  '

# Restore default (Normal) style so no stray formatting is introduced
$textRange.Style = "Normal"
